# Matrix dig wbs monitor - update for database (MySQL) integration
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hidden defined name used by the MySQL add-in/export tooling
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", '=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&" "&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)')
$definedName.Visible = $false

# Fill in the new matrix intersections with a lowercase "v"
$ws.Range("D2").Value = "v"
$ws.Range("B3").Value = "v"
$ws.Range("D3").Value = "v"
$ws.Range("B4").Value = "v"
$ws.Range("C4").Value = "v"

# Update the active selection to match the author's cursor position
$ws.Range("D3").Select()
